# Update cosinor-per-day statistics re-computed after re-running the
# CircaDB / CircadiPy analyses ("Make figures again to publication").
# Only the numeric/confidence-interval results for rows 2-10 change;
# labels, formatting and layout are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "[54.07623963002415, 72.37257851630358]"
$ws.Range("U2").Value = "[43.99843776378218, 56.5652688064105]"
$ws.Range("M3").Value = "[52.747887138107615, 73.32256263076529]"
$ws.Range("N3").Value = 0.0000000000000004440892098500626
$ws.Range("O3").Value = 0.0000000000000004440892098500626
$ws.Range("Q3").Value = "[1.2138686329185786, 1.5660792207084269]"
$ws.Range("U3").Value = "[42.17554580035242, 55.03304690542219]"
$ws.Range("Y3").Value = 18.61111111111143
$ws.Range("Z3").Value = 20.00074074074109
$ws.Range("M4").Value = "[51.74863700777345, 74.33454061030149]"
$ws.Range("N4").Value = 0.00000000000001154631945610163
$ws.Range("O4").Value = 0.00000000000001154631945610163
$ws.Range("Q4").Value = "[1.0755001877154244, 1.4528686746331188]"
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("U4").Value = "[42.33682188251351, 55.38733959332779]"
$ws.Range("Y4").Value = 19.05777777777811
$ws.Range("Z4").Value = 20.54666666666703
$ws.Range("M5").Value = "[50.539862095766026, 76.18154961532997]"
$ws.Range("N5").Value = 0.0000000000006012967901369848
$ws.Range("O5").Value = 0.0000000000006012967901369848
$ws.Range("Q5").Value = "[0.786184347745194, 1.1887107337907326]"
$ws.Range("R5").Value = 0.0000000000007540634783254063
$ws.Range("S5").Value = 0.0000000000007540634783254063
$ws.Range("U5").Value = "[42.70931548880267, 55.93753534121676]"
$ws.Range("Y5").Value = 20.10000000000035
$ws.Range("Z5").Value = 21.68814814814852
$ws.Range("M6").Value = "[52.1596358139696, 77.8783210049746]"
$ws.Range("N6").Value = 0.0000000000002930988785010413
$ws.Range("O6").Value = 0.0000000000002930988785010413
$ws.Range("Q6").Value = "[0.47171060864711656, 0.8993948938205012]"
$ws.Range("R6").Value = 0.0000000653271632344854
$ws.Range("S6").Value = 0.0000000653271632344854
$ws.Range("U6").Value = "[45.43053411967263, 58.90926377981312]"
$ws.Range("Y6").Value = 22.04692692692743
$ws.Range("Z6").Value = 23.79831831831886
$ws.Range("M7").Value = "[54.15592708301017, 77.48176839238522]"
$ws.Range("N7").Value = 0.000000000000008215650382226158
$ws.Range("O7").Value = 0.000000000000008215650382226158
$ws.Range("Q7").Value = "[0.09434212172942313, 0.4465527095192696]"
$ws.Range("R7").Value = 0.003398070425195687
$ws.Range("S7").Value = 0.003398070425195687
$ws.Range("U7").Value = "[43.857773451031626, 57.233182687041456]"
$ws.Range("Y7").Value = 23.90134134134188
$ws.Range("Z7").Value = 25.34366366366424
$ws.Range("M8").Value = "[53.821554040548534, 76.21641989639008]"
$ws.Range("N8").Value = 0.000000000000003108624468950438
$ws.Range("O8").Value = 0.000000000000003108624468950438
$ws.Range("U8").Value = "[45.42415442640336, 58.915539838647106]"
$ws.Range("M9").Value = "[49.58324392043407, 75.98537358614695]"
$ws.Range("N9").Value = 0.000000000001965316798191452
$ws.Range("O9").Value = 0.000000000001965316798191452
$ws.Range("U9").Value = "[41.99781288275247, 55.440456417304134]"
$ws.Range("M10").Value = "[50.88832741641615, 73.7361068478306]"
$ws.Range("N10").Value = 0.00000000000002531308496145357
$ws.Range("O10").Value = 0.00000000000002531308496145357
$ws.Range("Q10").Value = "[-1.2075791581366175, -0.8302106712189232]"
$ws.Range("U10").Value = "[40.82082616982602, 53.281735605463794]"
$ws.Range("Y10").Value = 3.399759759759833
$ws.Range("Z10").Value = 4.945105105105219
